$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.3447283333333333
$ws.Cells.Item(2, 8).Value = 1.034185
$ws.Cells.Item(2, 9).Value = 0.07368549602308437
$ws.Cells.Item(2, 10).Value = 0.07368549602308436
$ws.Cells.Item(2, 13).Value = 12.86269466666666
$ws.Cells.Item(2, 14).Value = 38.58808399999999
$ws.Cells.Item(2, 15).Value = 0.1337831063410017
$ws.Cells.Item(2, 16).Value = 0.1337831063410017
$ws.Cells.Item(2, 17).Value = 4.434135294615555
$ws.Cells.Item(2, 18).Value = 39.90721765153999
$ws.Cells.Item(2, 19).Value = 0.009857874550245754
$ws.Cells.Item(2, 20).Value = 0.009857874550245753
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.3447283333333333
$ws.Cells.Item(3, 8).Value = 1.034185
$ws.Cells.Item(3, 9).Value = 0.07368549602308437
$ws.Cells.Item(3, 10).Value = 0.07368549602308436
$ws.Cells.Item(3, 15).Value = 0.3593152390330854
$ws.Cells.Item(3, 16).Value = 0.3593152390330854
$ws.Cells.Item(3, 17).Value = 11.90921953350944
$ws.Cells.Item(3, 18).Value = 107.182975801585
$ws.Cells.Item(3, 19).Value = 0.02647632161680603
$ws.Cells.Item(3, 20).Value = 0.02647632161680602
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.3447283333333333
$ws.Cells.Item(4, 8).Value = 1.034185
$ws.Cells.Item(4, 9).Value = 0.07368549602308437
$ws.Cells.Item(4, 10).Value = 0.07368549602308436
$ws.Cells.Item(4, 13).Value = 18.65324433333334
$ws.Cells.Item(4, 14).Value = 55.95973300000001
$ws.Cells.Item(4, 15).Value = 0.1940098117012772
$ws.Cells.Item(4, 16).Value = 0.1940098117012772
$ws.Cells.Item(4, 17).Value = 6.430301830289445
$ws.Cells.Item(4, 18).Value = 57.872716472605
$ws.Cells.Item(4, 19).Value = 0.0142957092085538
$ws.Cells.Item(4, 20).Value = 0.0142957092085538
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.3447283333333333
$ws.Cells.Item(5, 8).Value = 1.034185
$ws.Cells.Item(5, 9).Value = 0.07368549602308437
$ws.Cells.Item(5, 10).Value = 0.07368549602308436
$ws.Cells.Item(5, 13).Value = 7.643308666666666
$ws.Cells.Item(5, 14).Value = 22.929926
$ws.Cells.Item(5, 15).Value = 0.07949699519803316
$ws.Cells.Item(5, 16).Value = 0.07949699519803316
$ws.Cells.Item(5, 17).Value = 2.634865057812222
$ws.Cells.Item(5, 18).Value = 23.71378552031
$ws.Cells.Item(5, 19).Value = 0.005857775523511829
$ws.Cells.Item(5, 20).Value = 0.005857775523511828
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.3447283333333333
$ws.Cells.Item(6, 8).Value = 1.034185
$ws.Cells.Item(6, 9).Value = 0.07368549602308437
$ws.Cells.Item(6, 10).Value = 0.07368549602308436
$ws.Cells.Item(6, 13).Value = 22.43995333333334
$ws.Cells.Item(6, 14).Value = 67.31986000000001
$ws.Cells.Item(6, 15).Value = 0.2333948477266026
$ws.Cells.Item(6, 16).Value = 0.2333948477266026
$ws.Cells.Item(6, 17).Value = 7.735687712677778
$ws.Cells.Item(6, 18).Value = 69.6211894141
$ws.Cells.Item(6, 19).Value = 0.01719781512396696
$ws.Cells.Item(6, 20).Value = 0.01719781512396696
$ws.Cells.Item(7, 9).Value = 0.9050707286284559
$ws.Cells.Item(7, 10).Value = 0.9050707286284558
$ws.Cells.Item(7, 13).Value = 12.86269466666666
$ws.Cells.Item(7, 14).Value = 38.58808399999999
$ws.Cells.Item(7, 15).Value = 0.1337831063410017
$ws.Cells.Item(7, 16).Value = 0.1337831063410017
$ws.Cells.Item(7, 17).Value = 54.4639892317151
$ws.Cells.Item(7, 18).Value = 490.1759030854359
$ws.Cells.Item(7, 19).Value = 0.1210831735342286
$ws.Cells.Item(7, 20).Value = 0.1210831735342286
$ws.Cells.Item(8, 9).Value = 0.9050707286284559
$ws.Cells.Item(8, 10).Value = 0.9050707286284558
$ws.Cells.Item(8, 15).Value = 0.3593152390330854
$ws.Cells.Item(8, 16).Value = 0.3593152390330854
$ws.Cells.Item(8, 19).Value = 0.3252057051989825
$ws.Cells.Item(8, 20).Value = 0.3252057051989824
$ws.Cells.Item(9, 9).Value = 0.9050707286284559
$ws.Cells.Item(9, 10).Value = 0.9050707286284558
$ws.Cells.Item(9, 13).Value = 18.65324433333334
$ws.Cells.Item(9, 14).Value = 55.95973300000001
$ws.Cells.Item(9, 15).Value = 0.1940098117012772
$ws.Cells.Item(9, 16).Value = 0.1940098117012772
$ws.Cells.Item(9, 17).Value = 78.9826801331119
$ws.Cells.Item(9, 18).Value = 710.8441211980071
$ws.Cells.Item(9, 19).Value = 0.1755926016375444
$ws.Cells.Item(9, 20).Value = 0.1755926016375444
$ws.Cells.Item(10, 9).Value = 0.9050707286284559
$ws.Cells.Item(10, 10).Value = 0.9050707286284558
$ws.Cells.Item(10, 13).Value = 7.643308666666666
$ws.Cells.Item(10, 14).Value = 22.929926
$ws.Cells.Item(10, 15).Value = 0.07949699519803316
$ws.Cells.Item(10, 16).Value = 0.07949699519803316
$ws.Cells.Item(10, 17).Value = 32.36375360715044
$ws.Cells.Item(10, 18).Value = 291.273782464354
$ws.Cells.Item(10, 19).Value = 0.07195040336765673
$ws.Cells.Item(10, 20).Value = 0.07195040336765672
$ws.Cells.Item(11, 9).Value = 0.9050707286284559
$ws.Cells.Item(11, 10).Value = 0.9050707286284558
$ws.Cells.Item(11, 13).Value = 22.43995333333334
$ws.Cells.Item(11, 14).Value = 67.31986000000001
$ws.Cells.Item(11, 15).Value = 0.2333948477266026
$ws.Cells.Item(11, 16).Value = 0.2333948477266026
$ws.Cells.Item(11, 17).Value = 95.01658932121556
$ws.Cells.Item(11, 18).Value = 855.14930389094
$ws.Cells.Item(11, 19).Value = 0.2112388448900437
$ws.Cells.Item(11, 20).Value = 0.2112388448900437
$ws.Cells.Item(12, 7).Value = 0.027522
$ws.Cells.Item(12, 8).Value = 0.082566
$ws.Cells.Item(12, 9).Value = 0.005882812712079546
$ws.Cells.Item(12, 10).Value = 0.005882812712079545
$ws.Cells.Item(12, 13).Value = 12.86269466666666
$ws.Cells.Item(12, 14).Value = 38.58808399999999
$ws.Cells.Item(12, 15).Value = 0.1337831063410017
$ws.Cells.Item(12, 16).Value = 0.1337831063410017
$ws.Cells.Item(12, 17).Value = 0.3540070826159999
$ws.Cells.Item(12, 18).Value = 3.186063743544
$ws.Cells.Item(12, 19).Value = 0.0007870209586443345
$ws.Cells.Item(12, 20).Value = 0.0007870209586443344
$ws.Cells.Item(13, 7).Value = 0.027522
$ws.Cells.Item(13, 8).Value = 0.082566
$ws.Cells.Item(13, 9).Value = 0.005882812712079546
$ws.Cells.Item(13, 10).Value = 0.005882812712079545
$ws.Cells.Item(13, 15).Value = 0.3593152390330854
$ws.Cells.Item(13, 16).Value = 0.3593152390330854
$ws.Cells.Item(13, 17).Value = 0.9507937361340001
$ws.Cells.Item(13, 18).Value = 8.557143625205999
$ws.Cells.Item(13, 19).Value = 0.002113784255827736
$ws.Cells.Item(13, 20).Value = 0.002113784255827735
$ws.Cells.Item(14, 7).Value = 0.027522
$ws.Cells.Item(14, 8).Value = 0.082566
$ws.Cells.Item(14, 9).Value = 0.005882812712079546
$ws.Cells.Item(14, 10).Value = 0.005882812712079545
$ws.Cells.Item(14, 13).Value = 18.65324433333334
$ws.Cells.Item(14, 14).Value = 55.95973300000001
$ws.Cells.Item(14, 15).Value = 0.1940098117012772
$ws.Cells.Item(14, 16).Value = 0.1940098117012772
$ws.Cells.Item(14, 17).Value = 0.5133745905420001
$ws.Cells.Item(14, 18).Value = 4.620371314878001
$ws.Cells.Item(14, 19).Value = 0.001141323386544432
$ws.Cells.Item(14, 20).Value = 0.001141323386544432
$ws.Cells.Item(15, 7).Value = 0.027522
$ws.Cells.Item(15, 8).Value = 0.082566
$ws.Cells.Item(15, 9).Value = 0.005882812712079546
$ws.Cells.Item(15, 10).Value = 0.005882812712079545
$ws.Cells.Item(15, 13).Value = 7.643308666666666
$ws.Cells.Item(15, 14).Value = 22.929926
$ws.Cells.Item(15, 15).Value = 0.07949699519803316
$ws.Cells.Item(15, 16).Value = 0.07949699519803316
$ws.Cells.Item(15, 17).Value = 0.210359141124
$ws.Cells.Item(15, 18).Value = 1.893232270116
$ws.Cells.Item(15, 19).Value = 0.0004676659339231161
$ws.Cells.Item(15, 20).Value = 0.000467665933923116
$ws.Cells.Item(16, 7).Value = 0.027522
$ws.Cells.Item(16, 8).Value = 0.082566
$ws.Cells.Item(16, 9).Value = 0.005882812712079546
$ws.Cells.Item(16, 10).Value = 0.005882812712079545
$ws.Cells.Item(16, 13).Value = 22.43995333333334
$ws.Cells.Item(16, 14).Value = 67.31986000000001
$ws.Cells.Item(16, 15).Value = 0.2333948477266026
$ws.Cells.Item(16, 16).Value = 0.2333948477266026
$ws.Cells.Item(16, 17).Value = 0.6175923956400001
$ws.Cells.Item(16, 18).Value = 5.55833156076
$ws.Cells.Item(16, 19).Value = 0.001373018177139927
$ws.Cells.Item(16, 20).Value = 0.001373018177139927
$ws.Cells.Item(17, 7).Value = 0.07186433333333334
$ws.Cells.Item(17, 8).Value = 0.215593
$ws.Cells.Item(17, 9).Value = 0.01536096263638017
$ws.Cells.Item(17, 10).Value = 0.01536096263638017
$ws.Cells.Item(17, 13).Value = 12.86269466666666
$ws.Cells.Item(17, 14).Value = 38.58808399999999
$ws.Cells.Item(17, 15).Value = 0.1337831063410017
$ws.Cells.Item(17, 16).Value = 0.1337831063410017
$ws.Cells.Item(17, 17).Value = 0.924368977090222
$ws.Cells.Item(17, 18).Value = 8.319320793811999
$ws.Cells.Item(17, 19).Value = 0.002055037297883003
$ws.Cells.Item(17, 20).Value = 0.002055037297883003
$ws.Cells.Item(18, 7).Value = 0.07186433333333334
$ws.Cells.Item(18, 8).Value = 0.215593
$ws.Cells.Item(18, 9).Value = 0.01536096263638017
$ws.Cells.Item(18, 10).Value = 0.01536096263638017
$ws.Cells.Item(18, 15).Value = 0.3593152390330854
$ws.Cells.Item(18, 16).Value = 0.3593152390330854
$ws.Cells.Item(18, 17).Value = 2.482674151034778
$ws.Cells.Item(18, 18).Value = 22.344067359313
$ws.Cells.Item(18, 19).Value = 0.005519427961469236
$ws.Cells.Item(18, 20).Value = 0.005519427961469236
$ws.Cells.Item(19, 7).Value = 0.07186433333333334
$ws.Cells.Item(19, 8).Value = 0.215593
$ws.Cells.Item(19, 9).Value = 0.01536096263638017
$ws.Cells.Item(19, 10).Value = 0.01536096263638017
$ws.Cells.Item(19, 13).Value = 18.65324433333334
$ws.Cells.Item(19, 14).Value = 55.95973300000001
$ws.Cells.Item(19, 15).Value = 0.1940098117012772
$ws.Cells.Item(19, 16).Value = 0.1940098117012772
$ws.Cells.Item(19, 17).Value = 1.340502968518778
$ws.Cells.Item(19, 18).Value = 12.064526716669
$ws.Cells.Item(19, 19).Value = 0.002980177468634471
$ws.Cells.Item(19, 20).Value = 0.002980177468634471
$ws.Cells.Item(20, 7).Value = 0.07186433333333334
$ws.Cells.Item(20, 8).Value = 0.215593
$ws.Cells.Item(20, 9).Value = 0.01536096263638017
$ws.Cells.Item(20, 10).Value = 0.01536096263638017
$ws.Cells.Item(20, 13).Value = 7.643308666666666
$ws.Cells.Item(20, 14).Value = 22.929926
$ws.Cells.Item(20, 15).Value = 0.07949699519803316
$ws.Cells.Item(20, 16).Value = 0.07949699519803316
$ws.Cells.Item(20, 17).Value = 0.5492812817908889
$ws.Cells.Item(20, 18).Value = 4.943531536118
$ws.Cells.Item(20, 19).Value = 0.001221150372941481
$ws.Cells.Item(20, 20).Value = 0.001221150372941481
$ws.Cells.Item(21, 7).Value = 0.07186433333333334
$ws.Cells.Item(21, 8).Value = 0.215593
$ws.Cells.Item(21, 9).Value = 0.01536096263638017
$ws.Cells.Item(21, 10).Value = 0.01536096263638017
$ws.Cells.Item(21, 13).Value = 22.43995333333334
$ws.Cells.Item(21, 14).Value = 67.31986000000001
$ws.Cells.Item(21, 15).Value = 0.2333948477266026
$ws.Cells.Item(21, 16).Value = 0.2333948477266026
$ws.Cells.Item(21, 17).Value = 1.612632286331111
$ws.Cells.Item(21, 18).Value = 14.51369057698
$ws.Cells.Item(21, 19).Value = 0.003585169535451982
$ws.Cells.Item(21, 20).Value = 0.003585169535451983
